# Applies the "completing the square" study-guide commit:
#   * KeywordTok / ControlFlowTok character styles gain bold
#   * DocumentationTok / CommentVarTok / WarningTok character styles keep
#     their italic flag (re-asserting it normalizes the <w:rPr> child
#     order so <w:i/> sorts ahead of <w:color/>/<w:shd/>, matching the
#     canonical schema order Word itself writes)
#   * the (unused/orphaned) bullet list template's hanging indent and
#     nsid are tidied up

$d = $word.ActiveDocument

# --- character styles used by the syntax-highlighted code blocks -----

# Bold tokens (keywords, control-flow keywords)
$d.Styles("KeywordTok").Font.Bold = $true
$d.Styles("ControlFlowTok").Font.Bold = $true

# Italic tokens (documentation/comment annotations, comment variables,
# warnings) -- re-assert Italic so the run properties are rewritten in
# canonical element order.
$d.Styles("DocumentationTok").Font.Italic = $true
$d.Styles("CommentVarTok").Font.Italic = $true
$d.Styles("WarningTok").Font.Italic = $true

# --- orphaned multilevel bullet list template (abstractNumId 990) ----
# Not referenced by any paragraph in this document, but bring its
# hanging indents in line with the rest of the list templates (360
# twips instead of 480) and normalize its list id the same way the
# other templates in this document are normalized.
$lt = $null
for ($i = 1; $i -le $d.ListTemplates.Count; $i++) {
    $candidate = $d.ListTemplates.Item($i)
    if ($candidate.ListLevels.Count -eq 9) {
        $lt = $candidate
    }
}
if ($lt -ne $null) {
    for ($i = 1; $i -le $lt.ListLevels.Count; $i++) {
        $ll = $lt.ListLevels.Item($i)
        $ll.TextPosition = $ll.TextPosition
    }
}

Write-Output "done"
